# Interdiff between v9 and v10
#
# 1) The textbox reading "[command is undoable]" is moved/resized and its
#    text is changed to "[addressbook is modified]" (typed as three runs:
#    "[" / "addressbook" / " is modified]").
# 2) The rounded-rectangle shape "Add addressbook to addressBookStateList
#    and clear redundant states" keeps the same visible text, but its
#    trailing " and clea" / "r redundant states" runs become one run
#    " and clear redundant states".

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shapes by their current text, rather than by a brittle
# index/name, so the script is resilient to shape reordering.
$undoableBox = $null
$addBox = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "[command is undoable]") {
            $undoableBox = $shp
        } elseif ($t.StartsWith("Add addressbook to addressBookStateList")) {
            $addBox = $shp
        }
    }
}

# ---- "[command is undoable]" -> "[addressbook is modified]" ----------
$tr = $undoableBox.TextFrame.TextRange
$tr.Text = "["
$tr.InsertAfter("addressbook") | Out-Null
$tr.InsertAfter(" is modified]") | Out-Null

$undoableBox.Left = 4885942 / $EMU_PER_POINT
$undoableBox.Top = 2161268 / $EMU_PER_POINT
$undoableBox.Width = 1472017 / $EMU_PER_POINT
$undoableBox.Height = 646587 / $EMU_PER_POINT

# ---- merge " and clea" + "r redundant states" into one run -----------
$addTr = $addBox.TextFrame.TextRange
$fullText = $addTr.Text
$newTail = " and clear redundant states"
$startIdx = $fullText.IndexOf(" and clea") + 1
$tailRange = $addTr.Characters($startIdx, $newTail.Length)
$tailRange.Text = $newTail
